$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulation result values (re-run with starting value 0 per commit message).
$ws.Cells.Item(4,3).Value = 0.709838121424159
$ws.Cells.Item(4,4).Value = -0.7043648548065328
$ws.Cells.Item(4,5).Value = 0.02251668865119556
$ws.Cells.Item(4,6).Value = -0.1048481457669956
$ws.Cells.Item(4,7).Value = -15.68933939568382
$ws.Cells.Item(4,10).Value = -1.008436362205496
$ws.Cells.Item(4,11).Value = 1.002846615199675
$ws.Cells.Item(4,12).Value = 0.0008386433658731404
$ws.Cells.Item(4,13).Value = -0.07396382173572996
$ws.Cells.Item(4,14).Value = -0.02346185170250373
$ws.Cells.Item(5,3).Value = 0.4951174013371049
$ws.Cells.Item(5,4).Value = -0.8688260828237517
$ws.Cells.Item(5,5).Value = 0.0792089098900229
$ws.Cells.Item(5,6).Value = 0.05452322129049693
$ws.Cells.Item(5,7).Value = -50.22228417872717
$ws.Cells.Item(5,10).Value = -0.7864008077140071
$ws.Cells.Item(5,11).Value = 0.6903147245254549
$ws.Cells.Item(5,12).Value = 0.06683330590921639
$ws.Cells.Item(5,13).Value = -0.2478746442946916
$ws.Cells.Item(5,14).Value = 0.003686079274810859
$ws.Cells.Item(6,3).Value = 0.7154367294421816
$ws.Cells.Item(6,4).Value = -0.6986775341180164
$ws.Cells.Item(6,5).Value = -0.2216962304415145
$ws.Cells.Item(6,6).Value = -0.7903510745430192
$ws.Cells.Item(6,7).Value = -2.565327538908753
$ws.Cells.Item(7,3).Value = 0.5741374184962693
$ws.Cells.Item(7,4).Value = -0.8187589615275936
$ws.Cells.Item(7,5).Value = 0.07926191665981962
$ws.Cells.Item(7,6).Value = -0.07491292752311646
$ws.Cells.Item(7,7).Value = 0.02250867104918054
$ws.Cells.Item(8,3).Value = -0.7044567187806103
$ws.Cells.Item(8,4).Value = 0.7097469514733179
$ws.Cells.Item(8,5).Value = 0.1262411707105019
$ws.Cells.Item(8,6).Value = -0.2431365207358937
$ws.Cells.Item(8,7).Value = 7.850925357428166
$ws.Cells.Item(9,3).Value = -0.4950578647948832
$ws.Cells.Item(9,4).Value = 0.8688600997818102
$ws.Cells.Item(9,5).Value = 0.07920438220063698
$ws.Cells.Item(9,6).Value = 0.05440984276865816
$ws.Cells.Item(9,7).Value = 45.50989226623803
$ws.Cells.Item(9,10).Value = 0.7897800162999066
$ws.Cells.Item(9,11).Value = -0.6935221643312095
$ws.Cells.Item(9,12).Value = 0.06669108186216435
$ws.Cells.Item(9,13).Value = -0.2490977430243676
$ws.Cells.Item(9,14).Value = 1.567155545411448
$ws.Cells.Item(10,10).Value = 5.481173967358281
$ws.Cells.Item(10,11).Value = 3.236682058183399
$ws.Cells.Item(10,12).Value = -1.072841450705498
$ws.Cells.Item(10,13).Value = 0.6003564245890757
$ws.Cells.Item(10,14).Value = -1.294525327934617
$ws.Cells.Item(11,3).Value = 0.7615661890233265
$ws.Cells.Item(11,4).Value = -0.6480871390026202
$ws.Cells.Item(11,5).Value = 0.06776645399750968
$ws.Cells.Item(11,6).Value = -0.1375952508607388
$ws.Cells.Item(11,7).Value = 1.54391159979668
$ws.Cells.Item(11,10).Value = -0.2537703572811663
$ws.Cells.Item(11,11).Value = 0.4663469400252133
$ws.Cells.Item(11,12).Value = 0.0765219550105843
$ws.Cells.Item(11,13).Value = -0.0901786962830054
$ws.Cells.Item(11,14).Value = 1.547280150959703
$ws.Cells.Item(12,3).Value = 0.9774860038907432
$ws.Cells.Item(12,4).Value = 0.2110002944970373
$ws.Cells.Item(12,5).Value = -0.7504605364110853
$ws.Cells.Item(12,6).Value = 0.7997849466677615
$ws.Cells.Item(12,7).Value = 2.079634326256628
$ws.Cells.Item(12,8).Value = -24.43232456647322
$ws.Cells.Item(12,10).Value = -0.8878371774424391
$ws.Cells.Item(12,11).Value = 0.8828946577927232
$ws.Cells.Item(12,12).Value = 0.000820642913212596
$ws.Cells.Item(12,13).Value = -0.07599525702237536
$ws.Cells.Item(12,14).Value = -0.02356525107181082
$ws.Cells.Item(12,15).Value = 1.136082408956695
$ws.Cells.Item(13,3).Value = 0.7516053174122794
$ws.Cells.Item(13,4).Value = -0.6596133396558523
$ws.Cells.Item(13,5).Value = 0.06674291432958364
$ws.Cells.Item(13,6).Value = -0.2482501801507858
$ws.Cells.Item(13,7).Value = 0.003685923931887904
$ws.Cells.Item(13,8).Value = -1.047831593073773
$ws.Cells.Item(14,3).Value = 0.9880942860847126
$ws.Cells.Item(14,4).Value = 0.1538495458973727
$ws.Cells.Item(14,5).Value = 0.06080772441915611
$ws.Cells.Item(14,6).Value = -0.2203478628010135
$ws.Cells.Item(14,7).Value = 3.223946284529929
$ws.Cells.Item(14,8).Value = 0.02358346987128835
$ws.Cells.Item(14,10).Value = 0.03260245534916581
$ws.Cells.Item(14,11).Value = 0.00002127956400327884
$ws.Cells.Item(14,12).Value = 0.05921318600546775
$ws.Cells.Item(14,13).Value = -0.223965678708776
$ws.Cells.Item(14,14).Value = -0.08308531615432913
$ws.Cells.Item(14,15).Value = -0.8088804624166962
$ws.Cells.Item(15,3).Value = 0.7033050888063973
$ws.Cells.Item(15,4).Value = 0.7108885709161049
$ws.Cells.Item(15,5).Value = 0.07375441594451834
$ws.Cells.Item(15,6).Value = -0.1054360420123764
$ws.Cells.Item(15,7).Value = 0.02475336635879038
$ws.Cells.Item(15,8).Value = -0.1271352062475561
$ws.Cells.Item(16,3).Value = -0.7090749030746691
$ws.Cells.Item(16,4).Value = 0.7051331686116852
$ws.Cells.Item(16,5).Value = 0.0008454985959567293
$ws.Cells.Item(16,6).Value = -0.07609126113049207
$ws.Cells.Item(16,7).Value = -1.594326038658507
$ws.Cells.Item(16,8).Value = 1.422240593000021
$ws.Cells.Item(16,10).Value = -380.9848901678577
$ws.Cells.Item(16,11).Value = 378.9142337246403
$ws.Cells.Item(16,12).Value = 0.0009084156548920764
$ws.Cells.Item(16,13).Value = -0.07591156862764555
$ws.Cells.Item(16,14).Value = 1.59387768663573
$ws.Cells.Item(16,15).Value = -0.002646912214969954
$ws.Cells.Item(17,3).Value = 0.7516505208070711
$ws.Cells.Item(17,4).Value = -0.6595617685997213
$ws.Cells.Item(17,5).Value = 0.06676183667599031
$ws.Cells.Item(17,6).Value = -0.2488013281547434
$ws.Cells.Item(17,7).Value = 1.567108715126415
$ws.Cells.Item(17,8).Value = 1.048794957372736
$ws.Cells.Item(17,10).Value = -2.502599357115216
$ws.Cells.Item(17,11).Value = 2.194255690912645
$ws.Cells.Item(17,12).Value = 0.06669520612804566
$ws.Cells.Item(17,13).Value = -0.2498797157698865
$ws.Cells.Item(17,14).Value = 1.567050643247503
$ws.Cells.Item(17,15).Value = -0.3145896719217678
$ws.Cells.Item(18,3).Value = -0.04541780833798722
$ws.Cells.Item(18,4).Value = -0.9989682584138067
$ws.Cells.Item(18,5).Value = 0.07565941105184355
$ws.Cells.Item(18,6).Value = -0.1908466566018453
$ws.Cells.Item(18,7).Value = -1.613383092264097
$ws.Cells.Item(18,8).Value = 0.01794490614628489
$ws.Cells.Item(18,10).Value = 1.610638435140205
$ws.Cells.Item(18,11).Value = 0.003877560919212601
$ws.Cells.Item(18,12).Value = 0.05867533818450962
$ws.Cells.Item(18,13).Value = -0.2236385236731886
$ws.Cells.Item(18,14).Value = 1.654191206751706
$ws.Cells.Item(18,15).Value = 0.01640883014325819
$ws.Cells.Item(19,3).Value = 0.7578409545179549
$ws.Cells.Item(19,4).Value = 0.6524393788035275
$ws.Cells.Item(19,5).Value = 0.07388718327957092
$ws.Cells.Item(19,6).Value = -0.1070429185261272
$ws.Cells.Item(19,7).Value = 1.546118945238499
$ws.Cells.Item(19,8).Value = 0.126365514847335
$ws.Cells.Item(19,10).Value = 0.2694245874272864
$ws.Cells.Item(19,11).Value = 0.3232423478894169
$ws.Cells.Item(19,12).Value = 0.07354145837067667
$ws.Cells.Item(19,13).Value = -0.1055495801041808
$ws.Cells.Item(19,14).Value = 1.546091192454078
$ws.Cells.Item(19,15).Value = 0.3053416392146496
$ws.Cells.Item(22,3).Value = 0.7000936373123587
$ws.Cells.Item(22,4).Value = 0.7140510478913726
$ws.Cells.Item(22,5).Value = 0.1021027043340019
$ws.Cells.Item(22,6).Value = -0.1892458303693788
$ws.Cells.Item(22,7).Value = 4.725730185531894
$ws.Cells.Item(22,8).Value = 3.478380489162676
$ws.Cells.Item(23,3).Value = 0.7071098409300689
$ws.Cells.Item(23,4).Value = 0.7071037214298725
$ws.Cells.Item(23,5).Value = 0.07852964957127639
$ws.Cells.Item(23,6).Value = -0.1616363979732369
$ws.Cells.Item(23,7).Value = 0.01125629533192462
$ws.Cells.Item(23,8).Value = 0.01725818096435809
$ws.Cells.Item(24,10).Value = -0.0001485196493705326
$ws.Cells.Item(24,11).Value = 0.00001381400287788994
$ws.Cells.Item(24,12).Value = 0.1020167588643707
$ws.Cells.Item(24,13).Value = -0.1892600273497442
$ws.Cells.Item(24,14).Value = 0.0007076745986418191
$ws.Cells.Item(24,15).Value = -0.0003621411749867451
$ws.Cells.Item(25,3).Value = -0.7315259168507862
$ws.Cells.Item(25,4).Value = 0.6818136877577514
$ws.Cells.Item(25,5).Value = 0.06796370088910662
$ws.Cells.Item(25,6).Value = -0.2253312913318439
$ws.Cells.Item(25,7).Value = -0.02016292546854981
$ws.Cells.Item(25,8).Value = -12.93796713673039
$ws.Cells.Item(27,3).Value = -0.9105388727003462
$ws.Cells.Item(27,4).Value = 0.4134234698939865
$ws.Cells.Item(27,5).Value = 0.07182844798550909
$ws.Cells.Item(27,6).Value = -0.108872348717634
$ws.Cells.Item(27,7).Value = -0.1035193405593236
$ws.Cells.Item(27,8).Value = -2.98058526336463
$ws.Cells.Item(28,3).Value = -0.1704167405842537
$ws.Cells.Item(28,4).Value = 0.9853720795457011
$ws.Cells.Item(28,5).Value = 0.09435243202051799
$ws.Cells.Item(28,6).Value = -0.09986668875507175
$ws.Cells.Item(28,7).Value = 0.08907233567849963
$ws.Cells.Item(28,8).Value = 0.02358971182780927
$ws.Cells.Item(28,9).Value = -0.001105000015157252
$ws.Cells.Item(28,10).Value = -0.8791880709594413
$ws.Cells.Item(28,11).Value = 0.8743520840195337
$ws.Cells.Item(28,12).Value = 0.0005417905300352542
$ws.Cells.Item(28,13).Value = -0.0764591037305877
$ws.Cells.Item(28,14).Value = -0.02313859286838172
$ws.Cells.Item(28,15).Value = 1.141644621649656
$ws.Cells.Item(28,16).Value = -0.0007415277058695335
$ws.Cells.Item(29,3).Value = 0.7531460506550851
$ws.Cells.Item(29,4).Value = -0.6578535413509916
$ws.Cells.Item(29,5).Value = 0.06666755875355909
$ws.Cells.Item(29,6).Value = -0.2491756862722024
$ws.Cells.Item(29,7).Value = 0.004124049604997121
$ws.Cells.Item(29,8).Value = -1.07555274140744
$ws.Cells.Item(29,9).Value = -2.620523819797156
$ws.Cells.Item(30,3).Value = -0.05395292895189671
$ws.Cells.Item(30,4).Value = 0.9985434800063084
$ws.Cells.Item(30,5).Value = 0.05105830091872841
$ws.Cells.Item(30,6).Value = -0.1954875353090587
$ws.Cells.Item(30,7).Value = -0.1877447440663715
$ws.Cells.Item(30,8).Value = -0.1286220179767376
$ws.Cells.Item(30,9).Value = -0.01966590587380392
$ws.Cells.Item(30,10).Value = 0.02956458681318301
$ws.Cells.Item(30,11).Value = 0.00008898982777003431
$ws.Cells.Item(30,12).Value = 0.05921688485624587
$ws.Cells.Item(30,13).Value = -0.2237725396091237
$ws.Cells.Item(30,14).Value = -0.07835397738370777
$ws.Cells.Item(30,15).Value = -0.8049020769029773
$ws.Cells.Item(30,16).Value = 0.3825401462632464
$ws.Cells.Item(31,3).Value = -0.6399409290937439
$ws.Cells.Item(31,4).Value = 0.7684241074777763
$ws.Cells.Item(31,5).Value = 0.02854685782067485
$ws.Cells.Item(31,6).Value = 0.2297228834772579
$ws.Cells.Item(31,7).Value = 0.03002158436976246
$ws.Cells.Item(31,8).Value = -0.1311001623395306
$ws.Cells.Item(31,9).Value = -1093.146085553029
$ws.Cells.Item(31,10).Value = 0.0005559026595429201
$ws.Cells.Item(31,11).Value = 0.0005781542851671488
$ws.Cells.Item(31,12).Value = 0.07853220704964289
$ws.Cells.Item(31,13).Value = -0.1616284827868015
$ws.Cells.Item(31,14).Value = 0.01125176237556763
$ws.Cells.Item(31,15).Value = -0.00001736438764454146
$ws.Cells.Item(31,16).Value = 0.0009996479900988332
$ws.Cells.Item(32,3).Value = 0.994187667068993
$ws.Cells.Item(32,4).Value = 0.1076611539406717
$ws.Cells.Item(32,5).Value = 0.06243112669240598
$ws.Cells.Item(32,6).Value = -0.2142122347735561
$ws.Cells.Item(32,7).Value = -0.06752724724328267
$ws.Cells.Item(32,8).Value = -0.02001023668062291
$ws.Cells.Item(32,10).Value = -0.7592748373184847
$ws.Cells.Item(32,11).Value = 0.7550716968140087
$ws.Cells.Item(32,12).Value = 0.0005295381186969019
$ws.Cells.Item(32,13).Value = -0.07600710405390634
$ws.Cells.Item(32,14).Value = -0.02326626212531237
$ws.Cells.Item(32,15).Value = 1.322177650832585
$ws.Cells.Item(33,3).Value = 0.7516356472842874
$ws.Cells.Item(33,4).Value = -0.6595786126886428
$ws.Cells.Item(33,5).Value = 0.06670817977653899
$ws.Cells.Item(33,6).Value = -0.2483572814097438
$ws.Cells.Item(33,7).Value = 0.003695481821577288
$ws.Cells.Item(33,8).Value = -1.047375987935854
$ws.Cells.Item(34,3).Value = 0.9926423774005686
$ws.Cells.Item(34,4).Value = 0.1210830730920692
$ws.Cells.Item(34,5).Value = 0.06045923614224346
$ws.Cells.Item(34,6).Value = -0.2221193797542841
$ws.Cells.Item(34,7).Value = -0.08246946870803593
$ws.Cells.Item(34,8).Value = -0.02407687121513873
$ws.Cells.Item(34,10).Value = 0.02726112220237813
$ws.Cells.Item(34,11).Value = 0.0001776680990342017
$ws.Cells.Item(34,12).Value = 0.05898655754058675
$ws.Cells.Item(34,13).Value = -0.2239167652587719
$ws.Cells.Item(34,14).Value = -0.08342983042415342
$ws.Cells.Item(34,15).Value = -0.9662895120170831
$ws.Cells.Item(35,3).Value = 0.7008243332110968
$ws.Cells.Item(35,4).Value = 0.7133341157866682
$ws.Cells.Item(35,5).Value = 0.0735371656838327
$ws.Cells.Item(35,6).Value = -0.1066707909404621
$ws.Cells.Item(35,7).Value = 0.02467783536997837
$ws.Cells.Item(35,8).Value = -0.1268920052844351

Write-Host "Updated 259 cells with refreshed simulation results."
